$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 12124.75
$ws.Range("I26").Value = 5499.5
$ws.Range("K26").Value = 5499.5
$ws.Range("M26").Value = -5155.5
$ws.Range("H46").Value = 1118.25
$ws.Range("I46").Value = 1040.5
$ws.Range("K46").Value = 3121.5
$ws.Range("M46").Value = -3002.5
$ws.Range("H51").Value = 9537.875
$ws.Range("I51").Value = 8767.166999999999
$ws.Range("K51").Value = 8767.166999999999
$ws.Range("M51").Value = -8283.166999999999
$ws.Range("H60").Value = 1118.25
$ws.Range("I60").Value = 1040.5
$ws.Range("K60").Value = 3121.5
$ws.Range("M60").Value = -2637.5
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = $null
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = $null
$ws.Range("H138").Value = 2362.4167
$ws.Range("J138").Value = 2596.4
$ws.Range("L138").Value = 7789.200000000001
$ws.Range("N138").Value = -18069.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4905282
$ws.Range("I32").Value = 6099405
$ws.Range("K32").Value = 6099405
$ws.Range("M32").Value = -6099118
$ws.Range("H38").Value = 4098.8
$ws.Range("I38").Value = 998
$ws.Range("J38").Value = 8750
$ws.Range("K38").Value = 998
$ws.Range("L38").Value = 8750
$ws.Range("M38").Value = -531
$ws.Range("N38").Value = -9684

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1600.3158
$ws.Range("I80").Value = 443
$ws.Range("J80").Value = 1817.3125
$ws.Range("K80").Value = 443
$ws.Range("L80").Value = 1817.3125
$ws.Range("M80").Value = 555
$ws.Range("N80").Value = -3813.3125
$ws.Range("H83").Value = 1600.3158
$ws.Range("I83").Value = 443
$ws.Range("J83").Value = 1817.3125
$ws.Range("K83").Value = 2215
$ws.Range("L83").Value = 9086.5625
$ws.Range("M83").Value = 2777
$ws.Range("N83").Value = -19070.5625
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992
$ws.Range("H126").Value = 125000
$ws.Range("J126").Value = 125000
$ws.Range("L126").Value = 125000
$ws.Range("N126").Value = -134880
$ws.Range("H140").Value = 88000
$ws.Range("J140").Value = 88000
$ws.Range("L140").Value = 88000
$ws.Range("N140").Value = -98360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 15000
$ws.Range("J21").Value = 15000
$ws.Range("L21").Value = 15000
$ws.Range("N21").Value = -15470
$ws.Range("H26").Value = 10508.5
$ws.Range("I26").Value = 5020
$ws.Range("K26").Value = 5020
$ws.Range("M26").Value = -4733
$ws.Range("H58").Value = 14144.723
$ws.Range("I58").Value = 6204.3477
$ws.Range("K58").Value = 6204.3477
$ws.Range("M58").Value = -6001.3477
$ws.Range("H86").Value = 13475.3125
$ws.Range("I86").Value = 14043.286
$ws.Range("K86").Value = 14043.286
$ws.Range("M86").Value = -12920.286
$ws.Range("H89").Value = 13475.3125
$ws.Range("I89").Value = 14043.286
$ws.Range("K89").Value = 70216.42999999999
$ws.Range("M89").Value = -64600.42999999999
$ws.Range("H99").Value = 3529.9092
$ws.Range("J99").Value = 2599.25
$ws.Range("L99").Value = 2599.25
$ws.Range("N99").Value = -5595.25
$ws.Range("H126").Value = 3529.9092
$ws.Range("J126").Value = 2599.25
$ws.Range("L126").Value = 7797.75
$ws.Range("N126").Value = -12737.75
$ws.Range("H132").Value = 79189030
$ws.Range("I132").Value = 3517.375
$ws.Range("J132").Value = 205885860
$ws.Range("K132").Value = 10552.125
$ws.Range("L132").Value = 617657580
$ws.Range("M132").Value = -8022.125
$ws.Range("N132").Value = -617662640
$ws.Range("H134").Value = 27032674
$ws.Range("I134").Value = 2272
$ws.Range("J134").Value = 90922710
$ws.Range("K134").Value = 6816
$ws.Range("L134").Value = 272768130
$ws.Range("M134").Value = -4281
$ws.Range("N134").Value = -272773200
$ws.Range("H136").Value = 14144.723
$ws.Range("I136").Value = 6204.3477
$ws.Range("K136").Value = 18613.0431
$ws.Range("M136").Value = -16063.0431

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 216.5
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = $null
$ws.Range("H115").Value = 1467.3334
$ws.Range("I115").Value = 701
$ws.Range("K115").Value = 2103
$ws.Range("M115").Value = -928
$ws.Range("H131").Value = 1451.12
$ws.Range("I131").Value = 868.4286
$ws.Range("J131").Value = 1494.9785
$ws.Range("K131").Value = 2605.2858
$ws.Range("L131").Value = 4484.9355
$ws.Range("M131").Value = 2434.7142
$ws.Range("N131").Value = -14564.9355

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 166.25
$ws.Range("I2").Value = 176.22223
$ws.Range("K2").Value = 176.22223
$ws.Range("M2").Value = -63.22223
$ws.Range("H46").Value = 4000
$ws.Range("I46").Value = 4000
$ws.Range("K46").Value = 4000
$ws.Range("M46").Value = -3844
$ws.Range("H102").Value = 13029.6875
$ws.Range("I102").Value = 13029.6875
$ws.Range("K102").Value = 13029.6875
$ws.Range("M102").Value = -11407.6875
$ws.Range("H126").Value = 6135.931
$ws.Range("J126").Value = 5613.5625
$ws.Range("L126").Value = 16840.6875
$ws.Range("N126").Value = -21780.6875
$ws.Range("H132").Value = 661054.6
$ws.Range("I132").Value = 5442
$ws.Range("J132").Value = 1555071.9
$ws.Range("K132").Value = 16326
$ws.Range("L132").Value = 4665215.699999999
$ws.Range("M132").Value = -13796
$ws.Range("N132").Value = -4670275.699999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 5776.8
$ws.Range("I32").Value = 2161.3333
$ws.Range("K32").Value = 2161.3333
$ws.Range("M32").Value = -1844.3333
$ws.Range("H40").Value = 2666.9
$ws.Range("I40").Value = 2666.9
$ws.Range("K40").Value = 2666.9
$ws.Range("M40").Value = -2530.9
$ws.Range("H46").Value = 2910.7693
$ws.Range("J46").Value = 4406.4287
$ws.Range("L46").Value = 4406.4287
$ws.Range("N46").Value = -4782.4287
$ws.Range("H55").Value = 2063.6
$ws.Range("I55").Value = 1660.2727
$ws.Range("J55").Value = 2556.5557
$ws.Range("K55").Value = 1660.2727
$ws.Range("L55").Value = 2556.5557
$ws.Range("M55").Value = -1487.2727
$ws.Range("N55").Value = -2902.5557
$ws.Range("H68").Value = 10287.091
$ws.Range("J68").Value = 10397.143
$ws.Range("L68").Value = 10397.143
$ws.Range("N68").Value = -11895.143
$ws.Range("H71").Value = 10287.091
$ws.Range("J71").Value = 10397.143
$ws.Range("L71").Value = 51985.715
$ws.Range("N71").Value = -59473.715
$ws.Range("H100").Value = 3516.4167
$ws.Range("I100").Value = 2875
$ws.Range("J100").Value = 4157.8335
$ws.Range("K100").Value = 2875
$ws.Range("L100").Value = 4157.8335
$ws.Range("M100").Value = -2334
$ws.Range("N100").Value = -5239.8335
$ws.Range("H111").Value = 100387
$ws.Range("J111").Value = 100387
$ws.Range("L111").Value = 100387
$ws.Range("N111").Value = -108567
$ws.Range("H122").Value = 7687.3335
$ws.Range("I122").Value = 6499.6665
$ws.Range("K122").Value = 19498.9995
$ws.Range("M122").Value = -17048.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 9963.333000000001
$ws.Range("I61").Value = 9963.333000000001
$ws.Range("K61").Value = 9963.333000000001
$ws.Range("M61").Value = -9671.333000000001
$ws.Range("H81").Value = 1158
$ws.Range("I81").Value = 1072.5
$ws.Range("J81").Value = 1500
$ws.Range("K81").Value = 2145
$ws.Range("L81").Value = 3000
$ws.Range("M81").Value = -1084
$ws.Range("N81").Value = -5122
$ws.Range("H84").Value = 1158
$ws.Range("I84").Value = 1072.5
$ws.Range("J84").Value = 1500
$ws.Range("K84").Value = 10725
$ws.Range("L84").Value = 15000
$ws.Range("M84").Value = -5421
$ws.Range("N84").Value = -25608
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = $null
$ws.Range("H122").Value = 6935.5
$ws.Range("I122").Value = 4365.6665
$ws.Range("K122").Value = 13096.9995
$ws.Range("M122").Value = -10646.9995
$ws.Range("H124").Value = 45189.5
$ws.Range("J124").Value = 45189.5
$ws.Range("L124").Value = 45189.5
$ws.Range("N124").Value = -55009.5
$ws.Range("H132").Value = 250918.66
$ws.Range("I132").Value = 1562.0526
$ws.Range("K132").Value = 4686.1578
$ws.Range("M132").Value = -2156.1578
